## Daily attendance processing - 2025-11-05 08:54:36
##
## Adds 13 new PARASITOLOGY attendance scan records (captured 2025-11-05,
## session "1") to the "Attendance" sheet, extends the sheet's AutoFilter
## / _FilterDatabase range to cover them, and refreshes the matching
## rollup figures for those 13 students on the "Summary" sheet (risk
## status, missed-percentage, sessions-needed, total-attended,
## total-missed and PARASITOLOGY-attended counters).

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("Summary")
$attendance = $wb.Worksheets.Item("Attendance")

# ---------------------------------------------------------------------
# Helpers: write a value into a cell while forcing it to be stored as
# literal text (Excel would otherwise silently reinterpret strings such
# as "221506", "05/11/2025" or "6.9%" as a number/date/percentage).
#
#  - Set-TextValue          : target cell should end up with NO explicit
#                              style (plain, unstyled data cell, as on the
#                              "Attendance" sheet).
#  - Set-TextValueKeepStyle  : target cell already carries a meaningful
#                              style (border/alignment/number format,
#                              as in the "Summary" sheet's Percentage
#                              column) that must be preserved.
# ---------------------------------------------------------------------
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

function Set-TextValueKeepStyle($cell, $text) {
    $fmt = $cell.NumberFormat
    $align = $cell.HorizontalAlignment
    $border = $cell.Borders.Item(9).LineStyle
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = $fmt
    $cell.HorizontalAlignment = $align
    $cell.Borders.Item(9).LineStyle = $border
}

# ---------------------------------------------------------------------
# 1. Append the new scan records to the "Attendance" sheet (rows 259-271)
# ---------------------------------------------------------------------

$newRecords = @(
    @("221506", "فاطمه احمد اسماعيل الناجي", "10:43:22"),
    @("221938", "مانويلا ناكوتا مارينو لوكالى", "10:43:24"),
    @("221574", "عمر براء رجب", "10:43:26"),
    @("221863", "ميسم ايمن زيدان", "10:43:27"),
    @("221587", "ليم عثمان السر الشمباتي", "10:43:29"),
    @("221617", "محمد رياض جمال خشفه", "10:43:31"),
    @("221569", "هبه جعفر محمد شوكت", "10:43:33"),
    @("221294", "هاله يحى ابكر ابراهيم", "10:43:34"),
    @("221494", "حسن الصادق مصطفى الحاج", "10:43:35"),
    @("221631", "رغد الحاج حسين عبدالمتعال", "10:43:37"),
    @("222004", "احمد ايمن احمد بشير", "10:43:39"),
    @("221599", "سلمى عبد الرحمن عبيد موسى", "10:43:40"),
    @("221629", "ناصر عبدالحميد الحسيني", "10:43:42")
)

$startRow = 259
$row = $startRow
foreach ($rec in $newRecords) {
    $studentId = $rec[0]
    $name = $rec[1]
    $time = $rec[2]

    Set-TextValue $attendance.Cells.Item($row, 1) $studentId
    $attendance.Cells.Item($row, 2).Value = $name
    $attendance.Cells.Item($row, 3).Value = "Year 2"
    $attendance.Cells.Item($row, 4).Value = "C1"
    $attendance.Cells.Item($row, 5).Value = "$studentId@med.asu.edu.eg"
    $attendance.Cells.Item($row, 6).Value = "PARASITOLOGY"
    Set-TextValue $attendance.Cells.Item($row, 7) "1"
    $attendance.Cells.Item($row, 8).Value = "PARASITOLOGY"
    Set-TextValue $attendance.Cells.Item($row, 9) "05/11/2025"
    $attendance.Cells.Item($row, 10).Value = $time
    $attendance.Cells.Item($row, 11).Value = "C1"

    $row = $row + 1
}

$lastRow = $row - 1

# Refresh the AutoFilter range so it spans the newly appended rows too.
$attendance.AutoFilterMode = $false
$attendance.Range("A1:K$lastRow").AutoFilter() | Out-Null

# The hidden workbook-level "_FilterDatabase" name for this sheet records
# the filter range explicitly and needs to be pointed at the new extent.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Attendance!_FilterDatabase") {
        $n.RefersTo = "='Attendance'!`$A`$1:`$K`$$lastRow"
    }
}

# ---------------------------------------------------------------------
# 2. Refresh the rollup numbers on the "Summary" sheet for each of the
#    13 students who now have one additional PARASITOLOGY attendance.
# ---------------------------------------------------------------------

# Reference cells (existing, unaffected rows) whose "Status" fill colour
# represents each risk tier - used to copy the correct style onto the
# cells that change tier.
$noRiskFill = $summary.Cells.Item(3, 6).Interior.Color
$lowRiskFill = $summary.Cells.Item(2, 6).Interior.Color
$moderateRiskFill = $summary.Cells.Item(4, 6).Interior.Color

function Get-RiskFill($status) {
    if ($status -eq "No Risk") { return $noRiskFill }
    elseif ($status -eq "Low Risk") { return $lowRiskFill }
    else { return $moderateRiskFill }
}

# row => (status, percentText, sessionsNeeded, totalAttended, totalMissed)
$summaryUpdates = @{
    107 = @("Low Risk",  "6.9%",  20, 2, 2)
    134 = @("No Risk",   "13.8%", 18, 4, 0)
    138 = @("Low Risk",  "6.9%",  20, 2, 2)
    152 = @("No Risk",   "10.3%", 19, 3, 1)
    153 = @("No Risk",   "10.3%", 19, 3, 1)
    158 = @("Low Risk",  "3.4%",  21, 1, 3)
    160 = @("Low Risk",  "6.9%",  20, 2, 2)
    163 = @("Low Risk",  "3.4%",  21, 1, 3)
    165 = @("Low Risk",  "6.9%",  20, 2, 2)
    166 = @("Low Risk",  "6.9%",  20, 2, 2)
    201 = @("Low Risk",  "6.9%",  20, 2, 2)
    217 = @("No Risk",   "10.3%", 19, 3, 1)
    233 = @("No Risk",   "13.8%", 18, 4, 0)
}

foreach ($r in $summaryUpdates.Keys) {
    $vals = $summaryUpdates[$r]
    $status = $vals[0]
    $percent = $vals[1]
    $sessionsNeeded = $vals[2]
    $totalAttended = $vals[3]
    $totalMissed = $vals[4]

    $fCell = $summary.Cells.Item($r, 6)
    $fCell.Value = $status
    $fCell.Interior.Color = Get-RiskFill $status

    Set-TextValueKeepStyle $summary.Cells.Item($r, 7) $percent
    $summary.Cells.Item($r, 8).Value = $sessionsNeeded
    $summary.Cells.Item($r, 12).Value = $totalAttended
    $summary.Cells.Item($r, 13).Value = $totalMissed
    # AI column = Attended PARASITOLOGY (Total)
    $summary.Cells.Item($r, 35).Value = 1
}
